$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct a few existing values ---
$ws.Range("B201").Value = 237.9664
$ws.Range("B209").Value = 223.2874
$ws.Range("B216").Value = 195.4933

# --- New daily data for August 2021 (04-08-2021 .. 31-08-2021) ---
# Use "X" as a placeholder separator so Excel's automatic date recognition
# does not kick in while the values are typed in, then swap "X" back to "-"
# via Find/Replace (which respects the Text number format applied below).
$dates = @(
    "04X08X2021","05X08X2021","06X08X2021","07X08X2021","08X08X2021",
    "09X08X2021","10X08X2021","11X08X2021","12X08X2021","13X08X2021",
    "14X08X2021","15X08X2021","16X08X2021","17X08X2021","18X08X2021",
    "19X08X2021","20X08X2021","21X08X2021","22X08X2021","23X08X2021",
    "24X08X2021","25X08X2021","26X08X2021","27X08X2021","28X08X2021",
    "29X08X2021","30X08X2021","31X08X2021"
)

$values = @(
    189.0465,183.3585,176.7975,167.256,158.4805,
    162.8177,167.6666,168.7202,164.8898,166.3455,
    163.8154,160.9456,166.8518,163.7594,165.3503,
    161.6042,156.5656,156.4447,158.1017,159.5063,
    162.1367,164.2676,163.4635,161.6277,160.7181,
    157.3638,160.3258,163.602
)

$startRow = 217
$endRow = $startRow + $dates.Length - 1
$colA = $ws.Range("A$startRow`:A$endRow")

$colA.NumberFormat = "@"
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
[void]$colA.Replace("X", "-")
$colA.Style = "Normal"
